$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix A3: make it a real number instead of text "2"
$ws.Range("A3").Value = 2

# Add new row 4 data (ID is text-formatted, like "2356")
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2356"
$ws.Range("B4").Value = "Destination1"
$ws.Range("C4").Value = 12.8293764059038
$ws.Range("D4").Value = 80.12243270874025
$ws.Range("E4").Value = 100
